# Apply the JLLE_COVID_hoje.xlsx update:
#  1. Correct a few previously-entered values in row 39/40 (J39, N39, R39, R40).
#  2. Add a new column S ("...19") = MORTES_DIA / CASOS_ATIVOS (M / R) for every
#     existing data row, left blank when CASOS_ATIVOS is blank or zero.
#  3. Append two new daily rows (106 and 107) with their own S value/blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix previously recorded values on rows 39 and 40 -------------------
$ws.Cells.Item(39, 10).Value = 616   # J39 TESTES
$ws.Cells.Item(39, 14).Value = 578   # N39 CASOS_DIA (acumulado)
$ws.Cells.Item(39, 18).Value = 45    # R39 CASOS_ATIVOS
$ws.Cells.Item(40, 18).Value = 2     # R40 CASOS_ATIVOS

# --- 2. New header for column S --------------------------------------------
$ws.Range("S1").Value = "...19"

# --- 3. Fill column S for rows 2 through 105 --------------------------------
for ($r = 2; $r -le 105; $r++) {
    $m = $ws.Cells.Item($r, 13).Value2
    $rr = $ws.Cells.Item($r, 18).Value2
    if ($rr -eq $null -or $rr -eq 0) {
        $ws.Cells.Item($r, 19).Style = "Normal"
    } else {
        $ws.Cells.Item($r, 19).Value = $m / $rr
    }
}

# --- 4. Append new row 106 --------------------------------------------------
$ws.Cells.Item(106, 1).Value = 105
$ws.Cells.Item(106, 2).Value = 630
$ws.Cells.Item(106, 3).Value = 39
$ws.Cells.Item(106, 4).Value = 800
$ws.Cells.Item(106, 5).Value = 33
$ws.Cells.Item(106, 6).Value = 1502
$ws.Cells.Item(106, 7).Value = 26
$ws.Cells.Item(106, 8).Value = 13
$ws.Cells.Item(106, 9).Value = 39
$ws.Cells.Item(106, 10).Value = 13283
$ws.Cells.Item(106, 11).Value = 8315
$ws.Cells.Item(106, 12).Value = 3466
$ws.Cells.Item(106, 13).Value = 323
$ws.Cells.Item(106, 14).Value = 9817
$ws.Cells.Item(106, 15).Value = 1
$ws.Cells.Item(106, 16).Value = 130
$ws.Cells.Item(106, 17).Value = 839
$ws.Cells.Item(106, 18).Value = 1388

# --- 5. Append new row 107 --------------------------------------------------
$ws.Cells.Item(107, 1).Value = 106
$ws.Cells.Item(107, 2).Value = 726
$ws.Cells.Item(107, 3).Value = 40
$ws.Cells.Item(107, 4).Value = 899
$ws.Cells.Item(107, 5).Value = 34
$ws.Cells.Item(107, 6).Value = 1699
$ws.Cells.Item(107, 7).Value = 26
$ws.Cells.Item(107, 8).Value = 14
$ws.Cells.Item(107, 9).Value = 40
$ws.Cells.Item(107, 10).Value = 14145
$ws.Cells.Item(107, 11).Value = 8879
$ws.Cells.Item(107, 12).Value = 3567
$ws.Cells.Item(107, 13).Value = 197
$ws.Cells.Item(107, 14).Value = 10578
$ws.Cells.Item(107, 15).Value = 1
$ws.Cells.Item(107, 16).Value = 96
$ws.Cells.Item(107, 17).Value = 939
$ws.Cells.Item(107, 18).Value = 761

# --- 6. Column S for the two new rows ---------------------------------------
for ($r = 106; $r -le 107; $r++) {
    $m = $ws.Cells.Item($r, 13).Value2
    $rr = $ws.Cells.Item($r, 18).Value2
    if ($rr -eq $null -or $rr -eq 0) {
        $ws.Cells.Item($r, 19).Style = "Normal"
    } else {
        $ws.Cells.Item($r, 19).Value = $m / $rr
    }
}
# Row 107 keeps S blank in the source workbook (not yet computed upstream).
$ws.Cells.Item(107, 19).ClearContents()
$ws.Cells.Item(107, 19).Style = "Normal"
